$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# Three partial transcripts of the Berkej manuscripts (rows 2, 3, 6) are now
# finished ("koncano") instead of "v delu", and get a green highlight.
$ws.Range("D2").Value = "končano"
$ws.Range("D3").Value = "končano"
$ws.Range("D6").Value = "končano"

# Row 5 (Akos Doncec / Martjanska II) status moves from "?" to "v delu".
$ws.Range("D5").Value = "v delu"

# --- Append new transcription entries ---
# Row 8: Maja Lampret - UKM Ms 139
$ws.Range("A8").Value = "Maja Lampret"
$ws.Range("B8").Value = "UKM Ms 139"
$ws.Range("C8").Value = "s. 84"
$ws.Range("D8").Value = "v delu"
$ws.Range("E8").Value = "da"

# Row 9: Maja Lampret - UKM Ms 598 (finished)
$ws.Range("A9").Value = "Maja Lampret"
$ws.Range("B9").Value = "UKM Ms 598"
$ws.Range("C9").Value = "s. 20"
$ws.Range("D9").Value = "končano"
$ws.Range("E9").Value = "da"

# Row 10: Spela Kovacic - NUK MS 1485-M13 (finished additional part)
$ws.Range("A10").Value = "Špela Kovačič"
$ws.Range("B10").Value = "NUK MS 1485-M13"
$ws.Range("C10").Value = "f. 145-191"
$ws.Range("D10").Value = "končano"
$ws.Range("E10").Value = "da"

# --- Highlight every "koncano" cell with the same green (Accent 6) fill ---
# Set the theme color on the first cell, then propagate the format (not the
# value) to the rest so only a single fill gets registered in the style table.
$ws.Range("D2").Interior.ThemeColor = 10
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Move selection to E2, matching the saved workbook's cursor position ---
$null = $ws.Range("E2").Select()
